$d = $word.ActiveDocument

# Locate the paragraph that holds the astromap credit/link line
# (the one crediting Jeník Hollan / CzechGlobe with the GaNight map link).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*CzechGlobe*GaNight*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $newText = "Jeník Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/"

    # Replace the whole (heavily run-fragmented / hyperlink-styled) paragraph
    # text with a single plain run, preserving the paragraph mark (and thus
    # the paragraph's own pPr/formatting).
    $r = $target.Range
    $r.End = $r.End - 1
    $r.Text = ""
    $r.InsertAfter($newText)

    # Re-insert a leading empty run ahead of the new text run, matching the
    # target structure (<w:r/><w:r><w:t>...</w:t></w:r>).
    $start = $target.Range.Start
    $rStart = $d.Range($start, $start)
    $emptyRunXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rStart.InsertXML($emptyRunXml)
}
